# Append the 2025 year-end poker figures (rows 201-210) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @(2025, "Andy",     1, 12, 0, 12, 44600, 60, 40,   "Active", 349),
    @(2025, "Prashant", 2,  9, 0,  9, 27800, 20,  0,   "Active", 365),
    @(2025, "Matt",     3,  9, 0,  9, 27050, 50, 30,   "Active", 362),
    @(2025, "Richard",  4,  6, 0,  6, 20800, 20,  0,   "Active", 366),
    @(2025, "Pepe",     5,  6, 0,  6, 20250,  0, -20,  "Active", 364),
    @(2025, "Maisy",    6,  5, 0,  5, 18350,  0, -20,  "Active", 360),
    @(2025, "Mark",     7,  5, 0,  5, 17050, 10,  0,   "Active", 361),
    @(2025, "Jon",      8,  4, 0,  4, 13050,  0, -20,  "Active", 357),
    @(2025, "Anthony",  9,  4, 0,  4,  9300, 10,  0,   "Active", 350),
    @(2025, "Alex",    10,  0, 0,  0,  4500,  0, -10,  "Active", 348)
)

$startRow = 201
for ($n = 0; $n -lt $data.Count; $n++) {
    $row = $startRow + $n
    $vals = $data[$n]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
    $ws.Cells.Item($row, 7).Value = $vals[6]
    $ws.Cells.Item($row, 8).Value = $vals[7]
    $ws.Cells.Item($row, 9).Value = $vals[8]
    $ws.Cells.Item($row, 10).Value = $vals[9]
    $ws.Cells.Item($row, 11).Value = $vals[10]
}

$ws.Range("A1:K210").Select() | Out-Null
